# Applies the para-page.docx edit described in the commit diff:
#   - third paragraph (the empty one right before the section break)
#     becomes centered and gets the repeated "Sample Text" run in
#     Times New Roman, 12pt
#   - the East Asian font used by the Normal / Heading styles (and the
#     document's run-property default) switches from "DejaVu Sans" to
#     "Tahoma"
#   - the List / Caption / Index styles pick up an explicit complex-script
#     (w:cs) font of "DejaVu Sans"

$d = $word.ActiveDocument

# --- 1. Third paragraph: center alignment + sample text run ------------
$targetPara = $d.Paragraphs.Item(3)
$targetPara.Format.Alignment = 1   # wdAlignParagraphCenter

$sampleText = "Sample Text Sample Text Sample Text Sample Text Sample Text " + `
              "Sample Text Sample Text Sample Text Sample Text Sample Text " + `
              "Sample Text Sample Text Sample Text Sample Text Sample Text " + `
              "Sample Text Sample Text Sample Text Sample Text Sample Text"

# Insert at a collapsed range so only the new run (not the paragraph
# mark's rPr) picks up the formatting.
$insertion = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
$insertion.InsertAfter($sampleText)

$textRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)
$textRange.Font.Name = "Times New Roman"
$textRange.Font.Size = 12

# --- 2. Swap the East Asian font on the styles that referenced it ------
$d.Styles("Normal").Font.NameFarEast = "Tahoma"
$d.Styles("Heading").Font.NameFarEast = "Tahoma"

# --- 3. Give List / Caption / Index an explicit complex-script font ----
$d.Styles("List").Font.NameBi = "DejaVu Sans"
$d.Styles("Caption").Font.NameBi = "DejaVu Sans"
$d.Styles("Index").Font.NameBi = "DejaVu Sans"
